$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column F (cell.Value column) to better fit numeric/short values
# (Excel COM stores ColumnWidth on a whole-pixel grid, so the closest
# representable width to the target 11.550625 characters is 11.5.)
$ws.Columns.Item(6).ColumnWidth = 10.666666666666666

# Update the "cell.Value" column (F) for rows 2-4 to hold the actual
# evaluated value of the corresponding formula in column C, instead of
# a copy of the formula text from column D.
# Row 2: C2 = A2+$B$2 = 1 + 2 = 3 (a number)
$ws.Range("F2").Value = 3
# Row 3: C3 = A3+$B$3 = 1 + 2 = 3 (a number)
$ws.Range("F3").Value = 3
# Row 4: C4 = "Test" & A4 & "R3C2" = "TestAR3C2" (a string)
$ws.Range("F4").Value = "TestAR3C2"
